$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (2-9) down to (3-10)
$ws.Rows.Item(2).Insert()

# Normalize formatting of the freshly-inserted row 2 to match the rest of the
# data rows: A column keeps the bordered/bold index style, B:E stay plain.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B3:E3").Copy()
$ws.Range("B2:E2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 2 data
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Tue, 07 Dec 2021 11:07:25 GMT"
$ws.Range("C2").Value = "Off the warpath: America 80 years after Pearl Harbour"
$ws.Range("D2").Value = "00:22:48"
$ws.Range("E2").Value = "https://sphinx.acast.com/theeconomistallaudio/theintelligencepodcast/offthewarpath-america80yearsafterpearlharbour/media.mp3"

# Renumber the index column (A) for the shifted rows 3-10 so it stays sequential 0-8
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

$wb.Save()
